# Add "canonical SMILES" (non-isomeric, stereo-stripped) as a new column D
# next to the existing "canonical isomeric SMILES" column C.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "canonical SMILES"

$ws.Range("D3").Value  = "c1ccc2c(c1)c(ncn2)Nc3cccc(c3)C(F)(F)F"
$ws.Range("D4").Value  = "c1ccc2c(c1)c(=Nc3cccc(c3)C(F)(F)F)[nH]cn2"
$ws.Range("D5").Value  = "c1ccc2c(c1)c(=[NH+]c3cccc(c3)C(F)(F)F)nc[nH]2"
$ws.Range("D6").Value  = "c1ccc2c(c1)c([nH+]cn2)Nc3cccc(c3)C(F)(F)F"
$ws.Range("D7").Value  = "c1ccc2c(c1)c(ncn2)[NH2+]c3cccc(c3)C(F)(F)F"
$ws.Range("D8").Value  = "c1ccc2c(c1)c(=Nc3cccc(c3)C(F)(F)F)nc[nH]2"
$ws.Range("D9").Value  = "c1ccc2c(c1)c(nc[nH+]2)[NH2+]c3cccc(c3)C(F)(F)F"
$ws.Range("D10").Value = "c1ccc2c(c1)c(ncn2)[N-]c3cccc(c3)C(F)(F)F"
$ws.Range("D11").Value = "c1ccc2c(c1)c([nH+]c[nH+]2)Nc3cccc(c3)C(F)(F)F"
$ws.Range("D12").Value = "c1ccc2c(c1)c([nH+]cn2)[NH2+]c3cccc(c3)C(F)(F)F"
$ws.Range("D13").Value = "c1ccc2c(c1)c([nH+]c[nH+]2)[NH2+]c3cccc(c3)C(F)(F)F"

# Target column width is 37.7109375 "raw" units, which is slightly off the
# character-width pixel grid that ColumnWidth rounds to (same as real Excel).
# 36.9 is the input that lands on the closest reachable grid point.
$ws.Columns.Item(4).ColumnWidth = 36.9
